$d = $word.ActiveDocument
$targetIndex = $d.Paragraphs.Count - 1
$p = $d.Paragraphs.Item($targetIndex)
$r0 = $p.Range
$r = $d.Range($r0.Start, $r0.Start)
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">Continuando tras la inicialización de la población creamos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>subproblemas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> asignándole a cada </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>subproblema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> el peso, una lista con los pesos de los vecinos más cercanos, el mejor individuo posible y su valor </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>fitness</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> (de la función </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>gte</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">Una vez preparados los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>subproblemas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> empezamos el bucle principal. Realizaremos el número de generaciones predefinidas en la entrada al programa -1 iteración ya que ya realizamos una evaluación previa en la inicialización y tenemos límite en el número de evaluaciones realizadas.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">El bucle itera sobre todos los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>subproblemas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> (N) y para cada </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>subproblema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> generaremos un número (G) de hijos usando la función de evolución descrita anteriormente y aprovechará el nuevo individuo para comprobar si es mejor que la mejor solución actual de todos los vecinos del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>subproblema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> (incluyéndose a sí mismo).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>La evaluación se realiza dependiendo de la función y de si tiene o no restricciones.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>Para terminar se realiza una gráfica con la población usando la librería ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>mathplotlib</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">’. Para comparar con mayor comodidad el resultado se aprovechan los archivos PF para pintar el frente ideal y posteriormente se pintan los objetivos conseguidos por cada </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>subproblema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>.</w:t></w:r></w:p>'
[void]$r.InsertXML($xmlFrag)
